# Update cryptocurrency price (column D) and 1h volume % (column E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "584.60", "0.310").
# Assign them with a leading apostrophe so Excel keeps the exact text
# (no trailing-zero / scientific-notation drift), then reset the cell
# style back to Normal so no stray text-format style is left on the cell.
$priceUpdates = @{
    'D2' = '65.943.17'
    'D3' = '3.015.14'
    'D5' = '584.60'
    'D6' = '161.55'
    'D7' = '0.999'
    'D8' = '3.011.21'
    'D9' = '0.518'
    'D10' = '6.84'
    'D12' = '0.460'
    'D13' = '0.0000253'
    'D16' = '65.921.78'
    'D17' = '3.516.86'
    'D19' = '3.015.59'
    'D20' = '457.35'
    'D21' = '13.96'
    'D23' = '7.39'
    'D24' = '82.38'
    'D26' = '12.42'
    'D27' = '10.67'
    'D31' = '0.0000106'
    'D33' = '27.18'
    'D36' = '0.995'
    'D37' = '5.84'
    'D38' = '2.16'
    'D39' = '49.91'
    'D40' = '2.97'
    'D41' = '0.310'
    'D43' = '43.43'
    'D44' = '8.45'
    'D45' = '390.59'
    'D46' = '2.798.83'
    'D47' = '0.0355'
    'D48' = '134.39'
    'D50' = '23.76'
}

foreach ($ref in $priceUpdates.Keys) {
    $ws.Range($ref).Value = "'" + $priceUpdates[$ref]
    $ws.Range($ref).Style = "Normal"
}

# Column E holds the 1h volume percentage text (already unambiguous text).
$volumeUpdates = @{
    'E3' = '  +3.91%  '
    'E6' = '  +12.69%  '
    'E7' = '  -0.20%  '
    'E8' = '  +3.86%  '
    'E9' = '  +3.77%  '
    'E10' = '  -2.37%  '
    'E11' = '  +6.18%  '
    'E12' = '  +7.47%  '
    'E13' = '  +9.07%  '
    'E14' = '  +7.09%  '
    'E15' = '  -0.45%  '
    'E16' = '  +6.63%  '
    'E17' = '  +3.94%  '
    'E18' = '  +7.28%  '
    'E19' = '  +4.02%  '
    'E20' = '  +6.63%  '
    'E21' = '  +7.71%  '
    'E22' = '  +6.07%  '
    'E23' = '  +7.82%  '
    'E24' = '  +4.81%  '
    'E25' = '  +12.52%  '
    'E26' = '  +3.33%  '
    'E27' = '  +4.74%  '
    'E28' = '  -0.04%  '
    'E29' = '  +15.58%  '
    'E30' = '  +17.20%  '
    'E31' = '  -6.24%  '
    'E32' = '  +4.09%  '
    'E33' = '  +6.52%  '
    'E34' = '  +4.46%  '
    'E35' = '  -0.04%  '
    'E36' = '  +4.24%  '
    'E37' = '  +8.86%  '
    'E38' = '  +14.32%  '
    'E39' = '  +2.24%  '
    'E40' = '  +2.97%  '
    'E41' = '  +16.74%  '
    'E42' = '  +7.23%  '
    'E43' = '  +5.54%  '
    'E44' = '  +3.77%  '
    'E45' = '  +12.71%  '
    'E46' = '  +3.51%  '
    'E47' = '  +5.80%  '
    'E48' = '  +1.30%  '
    'E49' = '  +0.00%  '
    'E50' = '  +10.22%  '
    'E51' = '  +4.58%  '
}

foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}
